# "Reunite xlsx and sysdata.rda" - refresh the medicare-tables.xlsx source
# data (sheet "indiv") so it matches the R package's rebuilt sysdata.rda.
# Only the raw input figures (lower thresholds / family thresholds / per-
# child add-ons) change; the ROUND(...) formulas in columns F and J are
# shared formulas that recompute automatically from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indiv")

# Row 4 (2015-16, single, not SAPTO): lower_threshold
$ws.Range("E4").Value = 21335

# Row 5 (2015-16, single, SAPTO): lower_threshold
$ws.Range("E5").Value = 33738

# Row 6 (2014-15, single, not SAPTO): lower_threshold,
# lower_family_threshold and lower_up_for_each_child
$ws.Range("E6").Value = 20896
$ws.Range("I6").Value = 35261
$ws.Range("K6").Value = 3238

# Row 7 (2014-15, single, SAPTO): lower_threshold and
# lower_up_for_each_child
$ws.Range("E7").Value = 33044
$ws.Range("K7").Value = 3238
